$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 527.8261
$ws.Range("I39").Value = 144.3077
$ws.Range("J39").Value = 1026.4
$ws.Range("K39").Value = 432.9231
$ws.Range("L39").Value = 3079.2
$ws.Range("M39").Value = -136.9231
$ws.Range("N39").Value = -3671.2

$ws.Range("H40").Value = 1801
$ws.Range("I40").Value = 1667.1818
$ws.Range("J40").Value = 2128.111
$ws.Range("K40").Value = 1667.1818
$ws.Range("L40").Value = 2128.111
$ws.Range("M40").Value = -1492.1818
$ws.Range("N40").Value = -2478.111

$ws.Range("H64").Value = 3145.4412
$ws.Range("I64").Value = 3107.7273
$ws.Range("J64").Value = 3163.4783
$ws.Range("K64").Value = 3107.7273
$ws.Range("L64").Value = 3163.4783
$ws.Range("M64").Value = -2859.7273
$ws.Range("N64").Value = -3659.4783

$ws.Range("H67").Value = 3145.4412
$ws.Range("I67").Value = 3107.7273
$ws.Range("J67").Value = 3163.4783
$ws.Range("K67").Value = 3107.7273
$ws.Range("L67").Value = 3163.4783
$ws.Range("M67").Value = -2249.7273
$ws.Range("N67").Value = -4879.478300000001

$ws.Range("H74").Value = 5736.1
$ws.Range("I74").Value = 6175.375
$ws.Range("J74").Value = 3979
$ws.Range("K74").Value = 6175.375
$ws.Range("L74").Value = 3979
$ws.Range("M74").Value = -5239.375
$ws.Range("N74").Value = -5851

$ws.Range("H77").Value = 5736.1
$ws.Range("I77").Value = 6175.375
$ws.Range("J77").Value = 3979
$ws.Range("K77").Value = 30876.875
$ws.Range("L77").Value = 19895
$ws.Range("M77").Value = -26196.875
$ws.Range("N77").Value = -29255

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 12564.934
$ws.Range("I141").Value = 14809.25
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 44427.75
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -39247.75
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1854276.5
$ws.Range("I63").Value = 2779166.5
$ws.Range("J63").Value = 4496.6665
$ws.Range("K63").Value = 2779166.5
$ws.Range("L63").Value = 4496.6665
$ws.Range("M63").Value = -2778480.5
$ws.Range("N63").Value = -5868.6665

$ws.Range("H66").Value = 1854276.5
$ws.Range("I66").Value = 2779166.5
$ws.Range("J66").Value = 4496.6665
$ws.Range("K66").Value = 13895832.5
$ws.Range("L66").Value = 22483.3325
$ws.Range("M66").Value = -13892400.5
$ws.Range("N66").Value = -29347.3325

$ws.Range("H135").Value = 36786.637
$ws.Range("J135").Value = 36786.637
$ws.Range("L135").Value = 36786.637
$ws.Range("N135").Value = -46926.637

$ws.Range("H139").Value = 35675
$ws.Range("J139").Value = 35675
$ws.Range("L139").Value = 35675
$ws.Range("N139").Value = -45955

$ws.Range("H140").Value = 18533
$ws.Range("J140").Value = 18533
$ws.Range("L140").Value = 18533
$ws.Range("N140").Value = -28893

$ws.Range("H141").Value = 18038.166
$ws.Range("J141").Value = 21045.8
$ws.Range("L141").Value = 21045.8
$ws.Range("N141").Value = -31405.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 11433.5
$ws.Range("I35").Value = 3067
$ws.Range("J35").Value = 19800
$ws.Range("K35").Value = 3067
$ws.Range("L35").Value = 19800
$ws.Range("M35").Value = -2757
$ws.Range("N35").Value = -20420

$ws.Range("H81").Value = 10744.167
$ws.Range("J81").Value = 10744.167
$ws.Range("L81").Value = 10744.167
$ws.Range("N81").Value = -12866.167

$ws.Range("H82").Value = 18622.75
$ws.Range("I82").Value = 9640
$ws.Range("J82").Value = 21617
$ws.Range("K82").Value = 9640
$ws.Range("L82").Value = 21617
$ws.Range("M82").Value = -9257
$ws.Range("N82").Value = -22383

$ws.Range("H84").Value = 10744.167
$ws.Range("J84").Value = 10744.167
$ws.Range("L84").Value = 32232.501
$ws.Range("N84").Value = -42840.501

$ws.Range("H85").Value = 18622.75
$ws.Range("I85").Value = 9640
$ws.Range("J85").Value = 21617
$ws.Range("K85").Value = 9640
$ws.Range("L85").Value = 21617
$ws.Range("M85").Value = -8314
$ws.Range("N85").Value = -24269

$ws.Range("H86").Value = 1698.8948
$ws.Range("I86").Value = 1614.8334
$ws.Range("J86").Value = 1843
$ws.Range("K86").Value = 1614.8334
$ws.Range("L86").Value = 1843
$ws.Range("M86").Value = -491.8334
$ws.Range("N86").Value = -4089

$ws.Range("H89").Value = 1698.8948
$ws.Range("I89").Value = 1614.8334
$ws.Range("J89").Value = 1843
$ws.Range("K89").Value = 8074.166999999999
$ws.Range("L89").Value = 9215
$ws.Range("M89").Value = -2458.166999999999
$ws.Range("N89").Value = -20447

$ws.Range("H138").Value = 33635.617
$ws.Range("J138").Value = 33635.617
$ws.Range("L138").Value = 33635.617
$ws.Range("N138").Value = -43915.617

$ws.Range("H140").Value = 27888.889
$ws.Range("J140").Value = 27888.889
$ws.Range("L140").Value = 27888.889
$ws.Range("N140").Value = -38248.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 845.6429000000001
$ws.Range("I22").Value = 1264.4286
$ws.Range("J22").Value = 426.85715
$ws.Range("K22").Value = 1264.4286
$ws.Range("L22").Value = 426.85715
$ws.Range("M22").Value = -914.4286
$ws.Range("N22").Value = -1126.85715

$ws.Range("H138").Value = 38180
$ws.Range("J138").Value = 38180
$ws.Range("L138").Value = 38180
$ws.Range("N138").Value = -48460

$ws.Range("H140").Value = 62199.6
$ws.Range("J140").Value = 62199.6
$ws.Range("L140").Value = 62199.6
$ws.Range("N140").Value = -72559.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 816.6900000000001
$ws.Range("J131").Value = 865.9888999999999
$ws.Range("L131").Value = 2597.9667
$ws.Range("N131").Value = -12677.9667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3669
$ws.Range("I5").Value = 420.75
$ws.Range("J5").Value = 8000
$ws.Range("K5").Value = 420.75
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = -308.75
$ws.Range("N5").Value = -8224

$ws.Range("H134").Value = 28288
$ws.Range("J134").Value = 28288
$ws.Range("L134").Value = 84864
$ws.Range("N134").Value = -89934

$ws.Range("H140").Value = 29000
$ws.Range("J140").Value = 29000
$ws.Range("L140").Value = 29000
$ws.Range("N140").Value = -39360

$ws.Range("H141").Value = 39900
$ws.Range("J141").Value = 39900
$ws.Range("L141").Value = 39900
$ws.Range("N141").Value = -50260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4456.7856
$ws.Range("I46").Value = 921.9524
$ws.Range("K46").Value = 921.9524
$ws.Range("M46").Value = -733.9524

$ws.Range("H138").Value = 60314.5
$ws.Range("J138").Value = 60314.5
$ws.Range("L138").Value = 60314.5
$ws.Range("N138").Value = -70594.5

$ws.Range("H140").Value = 42549.445
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 42549.445
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 42549.445
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -52909.445

$ws.Range("H141").Value = 54550.145
$ws.Range("J141").Value = 54550.145
$ws.Range("L141").Value = 54550.145
$ws.Range("N141").Value = -64910.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 32600.666
$ws.Range("J140").Value = 32600.666
$ws.Range("L140").Value = 32600.666
$ws.Range("N140").Value = -42960.666

$ws.Range("H141").Value = 38714.285
$ws.Range("J141").Value = 38714.285
$ws.Range("L141").Value = 38714.285
$ws.Range("N141").Value = -49074.285
